$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '69.717.37'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').Value = '3.528.38'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '609.50'
$ws.Range('E5').Value = '  +3.44%  '
$ws.Range('D6').Value = '183.93'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('D7').Value = '0.613'
$ws.Range('E7').Value = '  -1.60%  '
$ws.Range('E8').Value = '  -0.13%  '
$ws.Range('E9').Value = '  +4.63%  '
$ws.Range('D10').Value = '0.639'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').Value = '53.43'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').Value = '0.0000307'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').Value = '9.43'
$ws.Range('E13').Value = '  -1.35%  '
$ws.Range('D14').Value = '4.094.90'
$ws.Range('E14').Value = '  -1.28%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '69.823.10'
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('B16').Value = 'BitcoinCash'
$ws.Range('C16').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D16').Value = '585.68'
$ws.Range('E16').Value = '  +4.25%  '
$ws.Range('D17').Value = '3.566.67'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '12.59'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('D19').Value = '18.82'
$ws.Range('E19').Value = '  -4.44%  '
$ws.Range('D21').Value = '0.987'
$ws.Range('E21').Value = '  -2.96%  '
$ws.Range('D22').Value = '17.47'
$ws.Range('E22').Value = '  -2.22%  '
$ws.Range('D23').Value = '4.67'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = '4.83'
$ws.Range('E24').Value = '  -1.71%  '
$ws.Range('D25').Value = '96.25'
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('D26').Value = '2.96'
$ws.Range('E26').Value = '  -0.54%  '
$ws.Range('D27').Value = '10.95'
$ws.Range('E27').Value = '  -5.22%  '
$ws.Range('D28').Value = '9.50'
$ws.Range('E28').Value = '  +3.72%  '
$ws.Range('D29').Value = '31.94'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('D30').Value = '6.97'
$ws.Range('E30').Value = '  -4.06%  '
$ws.Range('D31').Value = '12.11'
$ws.Range('E31').Value = '  -3.15%  '
$ws.Range('D32').Value = '0.113'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').Value = '63.19'
$ws.Range('E33').Value = '  -3.05%  '
$ws.Range('E34').Value = '  -1.84%  '
$ws.Range('E35').Value = '  +17.24%  '
$ws.Range('D36').Value = '530.67'
$ws.Range('E36').Value = '  -6.09%  '
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = '36.99'
$ws.Range('E39').Value = '  -2.95%  '
$ws.Range('D40').Value = '3.526.40'
$ws.Range('E40').Value = '  +5.34%  '
$ws.Range('D41').Value = '0.0₃0773'
$ws.Range('E41').Value = '  -0.53%  '
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('D43').Value = '0.135'
$ws.Range('D44').Value = '0.0454'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').Value = '2.93'
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D46').Value = '3.37'
$ws.Range('E46').Value = '  -5.75%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.141'
$ws.Range('E47').Value = '  +2.71%  '
$ws.Range('D48').Value = '9.10'
$ws.Range('E48').Value = '  -2.82%  '
$ws.Range('E49').Value = '  +0.38%  '
$ws.Range('D50').Value = '1.41'
$ws.Range('E50').Value = '  -4.04%  '
$ws.Range('D51').Value = '134.44'
$ws.Range('E51').Value = '  -1.83%  '
